# City index update 25-10 (October 2025) — adds the new month's figures to
# "byindeks_aarlig" (year-over-year city index) and appends a new 12-month
# rolling window row to "by_glid_indeks".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "byindeks_aarlig": refresh the running year (2024-2025) row so it
# now reflects data through October instead of September, and update the
# prior full-year standard error that shifts slightly with the new point.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("byindeks_aarlig")

$ws3.Range("G2").Value = 0.9676465159429686

$ws3.Range("C3").Value = 10
$ws3.Range("E3").Value = 1.0204
$ws3.Range("G3").Value = 0.9687275452614585
$ws3.Range("K3").Value = "okt"
$ws3.Range("L3").Value = "jan-okt"

# ---------------------------------------------------------------------
# Sheet "by_glid_indeks": the trailing-12-month window ending September
# 2025 is recomputed very slightly, and a brand new window ending October
# 2025 is appended as row 12.
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("by_glid_indeks")

$ws4.Range("A11").Value = 0.9991532249526388
$ws4.Range("B11").Value = -0.08467750473611924
$ws4.Range("E11").Value = 5.216467590334918
$ws4.Range("F11").Value = 1.315691968381006

$ws4.Range("A12").Value = 1.001402751121869
$ws4.Range("B12").Value = 0.1402751121869494
$ws4.Range("C12").Value = 20
$ws4.Range("D12").Value = 15.71971180997949
$ws4.Range("E12").Value = 5.326593606543288
$ws4.Range("F12").Value = 1.343467836346448
$ws4.Range("G12").Value = -2.7
$ws4.Range("H12").Value = 3
$ws4.Range("I12").Value = "2023 - (nov 2024 - okt 2025)"
$ws4.Range("J12").Value = 45931
$ws4.Range("K12").Value = 10
$ws4.Range("L12").Value = 2025
$ws4.Range("M12").Value = "12_months"

# Make sure the date cell keeps/gets the yyyy-mm-dd format used by the rest
# of column J (it already inherits it from the column style, this is just
# a safety net in case that inheritance isn't picked up).
$ws4.Range("J12").NumberFormat = "yyyy-mm-dd"
